# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the second data row
# (aaaa355a-f478-4402-b8c3-9d02d8180801) across the Overview, zh-cn and de-de
# sheets, reflecting a new handback report run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-26 15:00:27"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-26 15:00:10"
$zhcn.Range("K3").Value = "2016-08-26 15:00:52"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-08-26 15:00:27"
$dede.Range("K3").Value = "2016-08-26 15:01:19"
